# The "Sheet" worksheet contained a bunch of now-unneeded helper cells
# (the T1/T2 column labels in row 1/2 and the numbered helper row 3) plus
# duplicated scheduling-result rows. Clean this up:
#   - Row 1: keep the original input list / result label cells (A1, B1),
#     updating their text, and drop the helper T1 labels (D1, G1).
#   - Row 2: keep the algorithm name (A2) and drop the helper T2 labels
#     (E2, F2, H2, I2).
#   - Row 3: delete entirely (it only held numbered helper cells 0..6),
#     which shifts the rows below it up by one.
#   - What were rows 5/6 (now 4/5) keep their content, with row 4's text
#     updated to match the new input list / result label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (the C3:I3 numbered helper row). This shifts row 5 -> 4 and
# row 6 -> 5 automatically, just like pressing "Delete Sheet Rows" in Excel.
$ws.Rows.Item(3).Delete()

# Remove the now-unused helper columns C:I (the T1/T2 labels in rows 1-2
# used to live here, along with the helper row that was just deleted).
$ws.Range("C1:I1").EntireColumn.Delete()

# Update the remaining text content.
$ws.Range("A1").Value = "2 2 5, 3 5 7, 1 5 10"
$ws.Range("B1").Value = "Unscheduable"
$ws.Range("A4").Value = "2 2 5, 3 5 7, 1 5 10"
$ws.Range("B4").Value = "Unscheduable"
